$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 78
$ws1.Range("F3").Value = 11841
$ws1.Range("F4").Value = 4
$ws1.Range("F6").Value = 352
$ws1.Range("F8").Value = 11760
$ws1.Range("F9").Value = 487
$ws1.Range("F10").Value = 1172
$ws1.Range("F12").Value = 50
$ws1.Range("F13").Value = 1775
$ws1.Range("F14").Value = 5824
$ws1.Range("F15").Value = 123

# Sheet "演出" (Shows)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 1

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 78
$ws4.Range("F5").Value = 11841
$ws4.Range("F6").Value = 4
$ws4.Range("F8").Value = 1
$ws4.Range("F9").Value = 352
$ws4.Range("F11").Value = 11760
$ws4.Range("F12").Value = 487
$ws4.Range("F13").Value = 1172
$ws4.Range("F15").Value = 50
$ws4.Range("F16").Value = 1775
$ws4.Range("F18").Value = 5824
$ws4.Range("F19").Value = 123
